# Add the new "Ireland" entry to the SFC Country Group mapping table.
# (Report regenerated for both 19437 and 60001, June 2021.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 39
$ws.Cells.Item($newRow, 1).Value = "IE"
$ws.Cells.Item($newRow, 2).Value = "Ireland"
$ws.Cells.Item($newRow, 3).Value = "Ireland"

# Mirror the author's final selection/view state (best effort).
$ws.Range("C" + $newRow).Select()
